# Re-sort the worksheet tabs: move "总计" (the summary sheet) so it becomes
# the first tab, ahead of "2020-Q4". This mirrors the commit's reordering of
# the <sheets> entries in xl/workbook.xml (总计 -> first position/sheetId 1,
# 2020-Q4 -> second position/sheetId 2), while each sheet keeps its own data
# untouched.

$wb = $excel.ActiveWorkbook

$totalSheet   = $wb.Worksheets.Item("总计")
$quarterSheet = $wb.Worksheets.Item("2020-Q4")

# Move "总计" in front of "2020-Q4" so the tab order becomes: 总计, 2020-Q4
$totalSheet.Move($quarterSheet)

# Worksheet object references can become stale/positionally-rebound once the
# tab order changes, so re-look-up the sheet by name before touching it again.
# Moving a sheet makes it the active one; restore "2020-Q4" as the
# selected/active sheet to match the original workbook's active tab.
$wb.Worksheets.Item("2020-Q4").Activate()
